# Apply the DSM restructuring edit described by the commit:
# "Packaged classes properly" -- TCPOctetStream.java and TCPState.java moved
# into a new `states\` folder (previously `concreteState\` / repo root), and
# the dependency matrix (DSM) values were refreshed to reflect the updated
# call/import relationships between the files.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: source file paths (rows 2-4 moved into the new `states` package) ---
$ws.Range("B2").Value = ".\states\TCPOctetStream.java"
$ws.Range("B3").Value = ".\context\TCPConnection.java"
$ws.Range("B4").Value = ".\states\TCPState.java"

# --- DSM matrix updates ---
# Row 3 (TCPConnection.java): now also depends on file (7)
$ws.Range("I3").Value = "Import,Call,Use"

# Row 4 (TCPState.java): dependency info on files (1) and (2) updated
$ws.Range("C4").Value = "Parameter"
$ws.Range("D4").Value = "Import,Call,Parameter:8"
$ws.Range("I4").Value = ""

# Row 5 (TCPListen.java): dependency info on files (2), (3) and (5) updated
$ws.Range("D5").Value = "Import,Call,Parameter"
$ws.Range("E5").Value = "Contain,Extend,Return,Import,Use:2"
$ws.Range("G5").Value = "Call,Use"

# Row 6 (TCPEstablished.java): dependency info on files (1),(2),(3),(4) updated
$ws.Range("C6").Value = "Import,Parameter"
$ws.Range("D6").Value = "Import,Call:2,Parameter:2"
$ws.Range("E6").Value = "Contain,Extend,Return,Import,Use:2"
$ws.Range("F6").Value = "Call,Use"

# Row 7 (Demo.java): relationship moved from column E to column D
$ws.Range("E7").Value = ""
$ws.Range("D7").Value = "Create,Import,Call:4,Contain"

# Row 8 (TCPClosed.java): dependency info on files (4) and (5) updated
$ws.Range("F8").Value = "Call,Use"
$ws.Range("G8").Value = "Call,Use"
